$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Description" text for gradeABC (B49)
$ws.Range("B49").Value = "Between 0 and 11. Sum of how many grades were at least C (>= 5)"

# Update the "Note" text for gradeABC (D49)
$ws.Range("D49").Value = "48 missings in the original dataset; 39 when merged with oxwaspLbw and 7 when merged also with oxwaspbp. Might be useful considering the ratio of A,B,C over the total number of grades obtained."

# Center the "Number of missings" value vertically
$ws.Range("C49").VerticalAlignment = -4108

# Wrap the Note text and let the row grow to fit
$ws.Range("D49").WrapText = $true
